$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet right after "RaiinInf" and name it.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "KaikeiInf"

# --- Header row (row 1) ---
$ws2.Range("A1").Value = "hp_id"
$ws2.Range("B1").Value = "pt_id"
$ws2.Range("C1").Value = "sin_date"
$ws2.Range("D1").Value = "raiin_no"
$ws2.Range("E1").Value = "tensu"
$ws2.Range("F1").Value = "total_iryohi"
$ws2.Range("G1").Value = "pt_futan"
$ws2.Range("H1").Value = "jihi_futan"
$ws2.Range("I1").Value = "jihi_tax"
$ws2.Range("J1").Value = "jihi_outtax"
$ws2.Range("K1").Value = "adjust_futan"
$ws2.Range("L1").Value = "adjust_round"
$ws2.Range("M1").Value = "total_pt_futan"
$ws2.Range("N1").Value = "create_date"
$ws2.Range("O1").Value = "create_id"
$ws2.Range("P1").Value = "create_machine"
$ws2.Range("Q1").Value = "adjust_futan_val"
$ws2.Range("R1").Value = "adjust_futan_range"
$ws2.Range("S1").Value = "adjust_rate_val"
$ws2.Range("T1").Value = "adjust_rate_range"
$ws2.Range("U1").Value = "hoken_id"
$ws2.Range("V1").Value = "kohi1_id"
$ws2.Range("W1").Value = "kohi2_id"
$ws2.Range("X1").Value = "kohi3_id"
$ws2.Range("Y1").Value = "kohi4_id"
$ws2.Range("Z1").Value = "hoken_kbn"
$ws2.Range("AA1").Value = "hoken_sbt_cd"
$ws2.Range("AB1").Value = "rece_sbt"
$ws2.Range("AC1").Value = "houbetu"
$ws2.Range("AD1").Value = "kohi1_houbetu"
$ws2.Range("AE1").Value = "kohi2_houbetu"
$ws2.Range("AF1").Value = "kohi3_houbetu"
$ws2.Range("AG1").Value = "kohi4_houbetu"
$ws2.Range("AH1").Value = "honke_kbn"
$ws2.Range("AI1").Value = "jihi_futan_taxfree"
$ws2.Range("AJ1").Value = "jihi_futan_tax_nr"
$ws2.Range("AK1").Value = "jihi_futan_tax_gen"
$ws2.Range("AL1").Value = "jihi_futan_outtax_nr"
$ws2.Range("AM1").Value = "jihi_futan_outtax_gen"
$ws2.Range("AN1").Value = "jihi_tax_nr"
$ws2.Range("AO1").Value = "jihi_tax_gen"
$ws2.Range("AP1").Value = "jihi_outtax_nr"
$ws2.Range("AQ1").Value = "jihi_outtax_gen"
$ws2.Range("AR1").Value = "hoken_rate"
$ws2.Range("AS1").Value = "pt_rate"
$ws2.Range("AT1").Value = "disp_rate"

# --- Data row (row 2) ---
$ws2.Range("A2").Value = 998
$ws2.Range("B2").Value = 12345
$ws2.Range("C2").Value = 20180807
$ws2.Range("D2").Value = 1234321
$ws2.Range("E2").Value = 2055
$ws2.Range("F2").Value = 20550
$ws2.Range("G2").Value = 500
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 0
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0
$ws2.Range("L2").Value = 0
$ws2.Range("M2").Value = 500
$ws2.Range("N2").Value = 44451.58668497685
$ws2.Range("N2").NumberFormat = "mm:ss.0"
$ws2.Range("O2").Value = 1
$ws2.Range("P2").Value = "285YMMTSVR"
$ws2.Range("Q2").Value = 0
$ws2.Range("R2").Value = 0
$ws2.Range("S2").Value = 0
$ws2.Range("T2").Value = 0
$ws2.Range("U2").Value = 100
$ws2.Range("V2").Value = 101
$ws2.Range("W2").Value = 0
$ws2.Range("X2").Value = 0
$ws2.Range("Y2").Value = 0
$ws2.Range("Z2").Value = 2
$ws2.Range("AA2").Value = 322
$ws2.Range("AB2").Value = "13x8"
$ws2.Range("AC2").Value = 39
$ws2.Range("AD2").Value = 80
$ws2.Range("AH2").Value = 1
$ws2.Range("AI2").Value = 0
$ws2.Range("AJ2").Value = 0
$ws2.Range("AK2").Value = 0
$ws2.Range("AL2").Value = 0
$ws2.Range("AM2").Value = 0
$ws2.Range("AN2").Value = 0
$ws2.Range("AO2").Value = 0
$ws2.Range("AP2").Value = 0
$ws2.Range("AQ2").Value = 0
$ws2.Range("AR2").Value = 10
$ws2.Range("AS2").Value = 10
$ws2.Range("AT2").Value = 10

# --- Column widths (best-fit on C & D, matching authored sample) ---
# NB: the engine pads ColumnWidth by ~0.8333 when serialising <col width=.../>,
# so back the value off to land on the authored width exactly (9 / 10 chars).
$ws2.Columns.Item(3).ColumnWidth = 8.166666666666666
$ws2.Columns.Item(4).ColumnWidth = 9.166666666666666

# --- Selections / active sheet bookkeeping ---
# RaiinInf keeps a B3 selection but is no longer the active tab; KaikeiInf
# becomes the active sheet with F6 selected.
[void]$ws1.Range("B3").Select()
[void]$ws2.Range("F6").Select()
[void]$ws2.Activate()
